# Updates the cryptos list (Sheet1) with new price/volume figures.
# Values are written as literal TEXT (matching the workbook's existing
# inline-string cell type) by prefixing the payload with a quote-prefix
# apostrophe, then resetting the cell style back to Normal so no stray
# number-format/style metadata is left behind.

function Set-TextValue {
    param($ws, $cellRef, $val)
    $ws.Range($cellRef).Formula = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "37.737.96"
Set-TextValue $ws "E2" "  +0.09%  "
Set-TextValue $ws "D3" "2.080.82"
Set-TextValue $ws "E3" "  -1.41%  "
Set-TextValue $ws "E4" "  +0.05%  "
Set-TextValue $ws "D5" "234.24"
Set-TextValue $ws "E5" "  -0.64%  "
Set-TextValue $ws "D6" "0.624"
Set-TextValue $ws "E6" "  -0.50%  "
Set-TextValue $ws "D7" "58.43"
Set-TextValue $ws "E7" "  +0.17%  "
Set-TextValue $ws "D9" "0.391"
Set-TextValue $ws "E9" "  -0.09%  "
Set-TextValue $ws "D10" "0.0785"
Set-TextValue $ws "E10" "  +0.23%  "
Set-TextValue $ws "E11" "  +2.69%  "
Set-TextValue $ws "D12" "14.92"
Set-TextValue $ws "E12" "  +2.49%  "
Set-TextValue $ws "D13" "2.386.33"
Set-TextValue $ws "E13" "  -1.37%  "
Set-TextValue $ws "D14" "21.01"
Set-TextValue $ws "E14" "  -1.34%  "
Set-TextValue $ws "D15" "0.773"
Set-TextValue $ws "E15" "  -2.10%  "
Set-TextValue $ws "E16" "  +1.58%  "
Set-TextValue $ws "D17" "2.077.67"
Set-TextValue $ws "E17" "  -1.38%  "
Set-TextValue $ws "D18" "37.629.25"
Set-TextValue $ws "E18" "  -0.05%  "
Set-TextValue $ws "D19" "6.18"
Set-TextValue $ws "E19" "  -0.27%  "
Set-TextValue $ws "D20" "71.15"
Set-TextValue $ws "E20" "  +1.32%  "
Set-TextValue $ws "D21" "0.0₃0833"
Set-TextValue $ws "E21" "  +0.99%  "
Set-TextValue $ws "D22" "228.26"
Set-TextValue $ws "E22" "  +0.42%  "
Set-TextValue $ws "D23" "0.999"
Set-TextValue $ws "D24" "2.40"
Set-TextValue $ws "E24" "  -0.97%  "
Set-TextValue $ws "E25" "  -1.12%  "
Set-TextValue $ws "D26" "168.99"
Set-TextValue $ws "E26" "  +0.32%  "
Set-TextValue $ws "E27" "  +2.91%  "
Set-TextValue $ws "D28" "8.97"
Set-TextValue $ws "E28" "  -0.28%  "
Set-TextValue $ws "B29" "ImmutableX"
Set-TextValue $ws "C29" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D29" "1.40"
Set-TextValue $ws "E29" "  -1.28%  "
Set-TextValue $ws "B30" "EthereumClassic"
Set-TextValue $ws "C30" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D30" "19.48"
Set-TextValue $ws "E30" "  +0.35%  "
Set-TextValue $ws "E31" "  +1.42%  "
Set-TextValue $ws "D32" "4.67"
Set-TextValue $ws "E32" "  +0.50%  "
Set-TextValue $ws "D33" "0.0631"
Set-TextValue $ws "E33" "  +1.36%  "
Set-TextValue $ws "D34" "4.64"
Set-TextValue $ws "E34" "  +1.27%  "
Set-TextValue $ws "D35" "2.48"
Set-TextValue $ws "E35" "  -3.98%  "
Set-TextValue $ws "E36" "  +2.54%  "
Set-TextValue $ws "D37" "3.39"
Set-TextValue $ws "E37" "  -4.02%  "
Set-TextValue $ws "E38" "  -0.07%  "
Set-TextValue $ws "D39" "5.38"
Set-TextValue $ws "E39" "  -4.41%  "
Set-TextValue $ws "D40" "0.0978"
Set-TextValue $ws "E40" "  +1.13%  "
Set-TextValue $ws "D41" "98.22"
Set-TextValue $ws "E41" "  +0.47%  "
Set-TextValue $ws "E42" "  +0.63%  "
Set-TextValue $ws "E43" "  -2.67%  "
Set-TextValue $ws "D44" "1.453.97"
Set-TextValue $ws "E44" "  -1.64%  "
Set-TextValue $ws "B45" "InjectiveProtocol"
Set-TextValue $ws "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D45" "16.64"
Set-TextValue $ws "E45" "  +5.97%  "
Set-TextValue $ws "B46" "FTXToken"
Set-TextValue $ws "C46" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D46" "4.31"
Set-TextValue $ws "E46" "  +3.06%  "
Set-TextValue $ws "E47" "  -0.87%  "
Set-TextValue $ws "D48" "1.06"
Set-TextValue $ws "E48" "  +0.50%  "
Set-TextValue $ws "D49" "7.40"
Set-TextValue $ws "E49" "  +0.98%  "
Set-TextValue $ws "E50" "  -0.41%  "
Set-TextValue $ws "D51" "2.270.23"
Set-TextValue $ws "E51" "  -1.54%  "
